# Apply cryptos.xlsx price/volume refresh (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.199.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "'1.802.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'314.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.5297"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.59%  "
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").Value = "'0.08010"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").Value = "'41.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "'6.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'20.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.335"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'1.803.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "'92.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'0.00001097"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("D19").Value = "'0.06610"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "'17.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "'5.976"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").Value = "'28.235.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "'11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'2.235"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'160.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").Value = "'20.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").Value = "'2.008.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "'2.379"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").Value = "'123.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("D31").Value = "'0.1090"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("D33").Value = "'3.665"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").Value = "'5.557"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").Value = "'0.07275"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Value = "'12.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.92%  "
$ws.Range("D37").Value = "'8.889"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "'0.2168"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").Value = "'0.02317"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'5.082"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").Value = "'0.6207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "'1.166"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").Value = "'1.372"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("D44").Value = "'13.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").Value = "'0.6007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("D46").Value = "'3.767"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").Value = "'126.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'1.213"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").Value = "'1.932"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("D50").Value = "'0.06834"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'73.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.26%  "
